$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.541.16"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.956.68"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'244.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("D7").Value = "'58.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.46%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("D10").Value = "'55.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.50%  "

$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").Value = "'22.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "'0.831"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "2.237.86"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").Value = "'13.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "'5.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").Value = "1.960.79"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "36.455.97"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Value = "'70.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").Value = "'230.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.92%  "

$ws.Range("D23").Value = "'5.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.83%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "'2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.14%  "

$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("D27").Value = "'9.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.44%  "

$ws.Range("D28").Value = "'162.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  +9.50%  "

$ws.Range("E30").Value = "  -1.46%  "

$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("E32").Value = "  +2.84%  "

$ws.Range("D33").Value = "'4.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.26%  "

$ws.Range("D34").Value = "'0.0634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.23%  "

$ws.Range("D35").Value = "'4.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.02%  "

$ws.Range("D36").Value = "'6.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.29%  "

$ws.Range("D39").Value = "'2.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.13%  "

$ws.Range("D40").Value = "'3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").Value = "'0.0988"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("E43").Value = "  -3.45%  "

$ws.Range("D44").Value = "'0.0211"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").Value = "'16.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("D46").Value = "1.371.57"
$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("E47").Value = "  -4.41%  "

$ws.Range("D48").Value = "'88.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.15%  "

$ws.Range("D49").Value = "'7.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.05%  "

$ws.Range("D50").Value = "'2.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").Value = "'46.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.37%  "
